$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("B1").Value = "Type"
$ws.Range("C1").Value = "Subject"
$ws.Range("D1").Value = "Validity From"
$ws.Range("F1").Value = "Expired Date"
$ws.Range("G1").Value = "Reminder Date"
$ws.Range("H1").Value = "Email"
$ws.Range("I1").Value = "Description"
$ws.Range("E1").Value = "Validity To"

# --- Row 2 ---
# Copy the date format already used on G2 onto the new D2:F2 cells
$ws.Range("G2").Copy()
$ws.Range("D2:F2").PasteSpecial(-4122)

# Move the note text from K2 to A2, keeping its formatting
$ws.Range("K2").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A2").Value = "*Note Type => General / Birthday"

# Remove the old K2 cell entirely
$ws.Range("K2").Clear()

# --- Column widths ---
$ws.Columns("A").ColumnWidth = 30

# --- Selection state ---
$ws.Range("B6").Select()
